$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.303000000000001
$ws.Range("E4").Value = 12.718

$ws.Range("E5").Value = 13.41

$ws.Range("B6").Value = 6.941
$ws.Range("E6").Value = 12.63

$ws.Range("B7").Value = 6.638

$ws.Range("B8").Value = 6.303
$ws.Range("E8").Value = 12.678

$ws.Range("B16").Value = 6.796000000000001
$ws.Range("E16").Value = 12.91

$ws.Range("B20").Value = 5.775999999999999

$ws.Range("B21").Value = 6.247

$ws.Range("E22").Value = 13.406
